$wb = $excel.ActiveWorkbook

$achieves = $wb.Worksheets.Item("Achieves")
$waves = $wb.Worksheets.Item("Waves")

# --- Update maxCondition (C) and reward (D) / maxAchieveLevel (E) columns ---
# Trim the trailing element off each comma-separated progression (5 values -> 4 values)
# and drop the maxCondition count from 5 to 4.

$achieves.Range("C2").Value = 4
$achieves.Range("D2").Value = "1000,2000,3000,4000"
$achieves.Range("E2").Value = "500,1000,1500,2000"

$achieves.Range("C3").Value = 4
$achieves.Range("D3").Value = "500,1000,1500,2000"
$achieves.Range("E3").Value = "500,1000,1500,2000"

$achieves.Range("C4").Value = 4
$achieves.Range("D4").Value = "100,200,300,400"
$achieves.Range("E4").Value = "500,1000,1500,2000"

$achieves.Range("C5").Value = 4
$achieves.Range("D5").Value = "100,300,500,700"
$achieves.Range("E5").Value = "500,1000,1500,2000"

$achieves.Range("C6").Value = 4
$achieves.Range("D6").Value = "5,10,15,20"
$achieves.Range("E6").Value = "500,1000,1500,2000"

$achieves.Range("C7").Value = 4
$achieves.Range("D7").Value = "4,8,12,16"
$achieves.Range("E7").Value = "500,1000,1500,2000"

$achieves.Range("C8").Value = 4
$achieves.Range("D8").Value = "3,6,9,12"
$achieves.Range("E8").Value = "500,1000,1500,2000"

$achieves.Range("C9").Value = 4
$achieves.Range("D9").Value = "2,5,8,11"
$achieves.Range("E9").Value = "500,1000,1500,2000"

# --- Sheet view / active tab changes ---
# Move the selected tab from "Waves" back to "Achieves", and update the
# selected cell on each sheet.
$waves.Range("B3").Select()

$achieves.Activate()
$achieves.Range("E11").Select()
